# Scheduled-runner refresh of the market-price / profit columns (H:N) on
# each leve-profit sheet. Only numeric value cells change; no formulas,
# formatting, or structure are touched.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1550
$ws.Range("I43").Value = 1600
$ws.Range("J43").Value = 1500
$ws.Range("K43").Value = 1600
$ws.Range("L43").Value = 1500
$ws.Range("M43").Value = -1531
$ws.Range("N43").Value = -1638
# Row 112
$ws.Range("H112").Value = 15626622
$ws.Range("I112").Value = 2123.75
$ws.Range("J112").Value = 31251122
$ws.Range("K112").Value = 6371.25
$ws.Range("L112").Value = 93753366
$ws.Range("M112").Value = -5263.25
$ws.Range("N112").Value = -93755582
# Row 116
$ws.Range("H116").Value = 4091.7778
$ws.Range("I116").Value = 3488.3333
$ws.Range("J116").Value = 4574.533
$ws.Range("K116").Value = 3488.3333
$ws.Range("L116").Value = 4574.533
$ws.Range("M116").Value = -46.33329999999978
$ws.Range("N116").Value = -11458.533

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 12624.667
$ws.Range("I37").Value = 2017
$ws.Range("J37").Value = 17928.5
$ws.Range("K37").Value = 2017
$ws.Range("L37").Value = 17928.5
$ws.Range("M37").Value = -1744
$ws.Range("N37").Value = -18474.5
# Row 44
$ws.Range("H44").Value = 17639.6
$ws.Range("J44").Value = 20424.5
$ws.Range("L44").Value = 20424.5
$ws.Range("N44").Value = -21400.5
# Row 55
$ws.Range("H55").Value = 20902
$ws.Range("J55").Value = 20902
$ws.Range("L55").Value = 20902
$ws.Range("N55").Value = -21532
# Row 61
$ws.Range("H61").Value = 2429.6296
$ws.Range("I61").Value = 1794.7368
$ws.Range("J61").Value = 3937.5
$ws.Range("K61").Value = 1794.7368
$ws.Range("L61").Value = 3937.5
$ws.Range("M61").Value = -1582.7368
$ws.Range("N61").Value = -4361.5
# Row 80
$ws.Range("H80").Value = 29306.334
$ws.Range("J80").Value = 29306.334
$ws.Range("L80").Value = 29306.334
$ws.Range("N80").Value = -31302.334
# Row 83
$ws.Range("H83").Value = 29306.334
$ws.Range("J83").Value = 29306.334
$ws.Range("L83").Value = 87919.00199999999
$ws.Range("N83").Value = -97903.00199999999
# Row 110
$ws.Range("H110").Value = 2475.3333
$ws.Range("I110").Value = 594.1667
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 594.1667
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = 1450.8333
$ws.Range("N110").Value = -14090
# Row 132
$ws.Range("H132").Value = 2972.282
$ws.Range("I132").Value = 2338.913
$ws.Range("J132").Value = 3882.75
$ws.Range("K132").Value = 7016.739
$ws.Range("L132").Value = 11648.25
$ws.Range("M132").Value = -4486.739
$ws.Range("N132").Value = -16708.25
# Row 136
$ws.Range("H136").Value = 2429.6296
$ws.Range("I136").Value = 1794.7368
$ws.Range("J136").Value = 3937.5
$ws.Range("K136").Value = 5384.2104
$ws.Range("L136").Value = 11812.5
$ws.Range("M136").Value = -2834.2104
$ws.Range("N136").Value = -16912.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2533.1538
$ws.Range("I134").Value = 2391.4
$ws.Range("J134").Value = 3444.4285
$ws.Range("K134").Value = 7174.200000000001
$ws.Range("L134").Value = 10333.2855
$ws.Range("M134").Value = -4639.200000000001
$ws.Range("N134").Value = -15403.2855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 38001
$ws.Range("I4").Value = 4999.5
$ws.Range("K4").Value = 4999.5
$ws.Range("M4").Value = -4887.5
# Row 22
$ws.Range("H22").Value = 2175.1667
$ws.Range("J22").Value = 3866.6667
$ws.Range("L22").Value = 3866.6667
$ws.Range("N22").Value = -4566.6667
# Row 131
$ws.Range("H131").Value = 26584.572
$ws.Range("J131").Value = 27837.23
$ws.Range("L131").Value = 27837.23
$ws.Range("N131").Value = -37917.23
# Row 132
$ws.Range("H132").Value = 3554.8823
$ws.Range("I132").Value = 3375.5454
$ws.Range("K132").Value = 10126.6362
$ws.Range("M132").Value = -7596.636200000001
# Row 141
$ws.Range("H141").Value = 28574.143
$ws.Range("I141").Value = 3731.6667
$ws.Range("J141").Value = 30903.125
$ws.Range("K141").Value = 3731.6667
$ws.Range("L141").Value = 30903.125
$ws.Range("M141").Value = 1448.3333
$ws.Range("N141").Value = -41263.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 1176.3334
$ws.Range("I109").Value = 412.2857
$ws.Range("J109").Value = 1844.875
$ws.Range("K109").Value = 1236.8571
$ws.Range("L109").Value = 5534.625
$ws.Range("M109").Value = -196.8571000000002
$ws.Range("N109").Value = -7614.625
# Row 113
$ws.Range("H113").Value = 2564920.5
$ws.Range("I113").Value = 14286230
$ws.Range("J113").Value = 884.09375
$ws.Range("K113").Value = 42858690
$ws.Range("L113").Value = 2652.28125
$ws.Range("M113").Value = -42856520
$ws.Range("N113").Value = -6992.28125
# Row 131
$ws.Range("H131").Value = 1464.0197
$ws.Range("I131").Value = 4714
$ws.Range("J131").Value = 1110.7609
$ws.Range("K131").Value = 14142
$ws.Range("L131").Value = 3332.2827
$ws.Range("M131").Value = -9102
$ws.Range("N131").Value = -13412.2827

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 14666.667
$ws.Range("J5").Value = 15125
$ws.Range("L5").Value = 15125
$ws.Range("N5").Value = -15349
# Row 33
$ws.Range("H33").Value = 82013.57000000001
$ws.Range("I33").Value = 4000
$ws.Range("J33").Value = 95015.836
$ws.Range("K33").Value = 4000
$ws.Range("L33").Value = 95015.836
$ws.Range("M33").Value = -3748
$ws.Range("N33").Value = -95519.836
# Row 132
$ws.Range("H132").Value = 3728.0193
$ws.Range("I132").Value = 3822.4517
$ws.Range("J132").Value = 3588.6191
$ws.Range("K132").Value = 11467.3551
$ws.Range("L132").Value = 10765.8573
$ws.Range("M132").Value = -8937.355100000001
$ws.Range("N132").Value = -15825.8573

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 5000
$ws.Range("J38").Value = 5000
$ws.Range("L38").Value = 5000
$ws.Range("N38").Value = -5820
# Row 46
$ws.Range("H46").Value = 1942.8572
$ws.Range("J46").Value = 3900
$ws.Range("L46").Value = 3900
$ws.Range("N46").Value = -4276
# Row 132
$ws.Range("H132").Value = 2832.9722
$ws.Range("I132").Value = 2358
$ws.Range("K132").Value = 7074
$ws.Range("M132").Value = -4544
# Row 136
$ws.Range("H136").Value = 4005739.2
$ws.Range("I136").Value = 7697430.5
$ws.Range("J136").Value = 6407.0835
$ws.Range("K136").Value = 23092291.5
$ws.Range("L136").Value = 19221.2505
$ws.Range("M136").Value = -23089741.5
$ws.Range("N136").Value = -24321.2505

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 1550.4
$ws.Range("J5").Value = 1550.4
$ws.Range("L5").Value = 1550.4
$ws.Range("N5").Value = -1774.4
# Row 122
$ws.Range("H122").Value = 1962.2
$ws.Range("I122").Value = 1652.75
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 4958.25
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -2508.25
$ws.Range("N122").Value = -14500
# Row 132
$ws.Range("H132").Value = 1764477.9
$ws.Range("I132").Value = 2328841.5
$ws.Range("K132").Value = 6986524.5
$ws.Range("M132").Value = -6983994.5

Write-Output "Applied all sheet updates"
